# Revert "Powerpoint writer: consolidate text run nodes."
#
# Splits the leading "Slide " / "an " run of each title / caption textbox
# into two separate runs ("Slide" + " ", "an" + " ") instead of one combined
# run, matching the pre-consolidation run layout. Re-assigning a
# Characters() sub-range's Text to itself forces the host to materialize it
# as its own <a:r> without touching formatting (rPr stays empty).

$p = $ppt.ActivePresentation

function Split-LeadingWord($shape, [int]$wordLen) {
    $tr = $shape.TextFrame.TextRange
    $word = $tr.Characters(1, $wordLen)
    $word.Text = $word.Text
}

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if (-not $sh.HasTextFrame) {
            continue
        }
        $text = $sh.TextFrame.TextRange.Text
        if ($text.StartsWith("Slide ")) {
            Split-LeadingWord $sh 5   # "Slide"
        } elseif ($text.StartsWith("an ")) {
            Split-LeadingWord $sh 2   # "an"
        }
    }
}
